$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the paragraph that currently reads:
#   "Open up the web app with command ... web browser" + bookmark + " )"
# ---------------------------------------------------------------------------
$target = $null
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Open up the web app*") {
        $target = $cand
        $targetIndex = $i
        break
    }
}

# 1) Re-write the text of that paragraph: split the single run in two so a
#    gramStart/gramEnd proofing pair wraps the word "browser )" - keep the
#    paragraph mark (and its rsid/pPr) untouched by only touching the inner
#    range.
$innerRange = $d.Range($target.Range.Start, $target.Range.End - 1)

$openUpBody = "<w:r><w:t xml:space=`"preserve`">Open up the web app with command " +
    [char]8211 +
    " ng serve -o (-o means open it will automatically open I on the web </w:t></w:r>" +
    "<w:proofErr w:type=`"gramStart`"/><w:r><w:t>browser )</w:t></w:r><w:proofErr w:type=`"gramEnd`"/>"

$openUpXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' + $openUpBody + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$innerRange.InsertXML($openUpXml)

# ---------------------------------------------------------------------------
# 2) Insert two brand-new ListParagraph bullets right after that paragraph:
#      - "Installing angular material for design -  npm install --save ..."
#      - an empty bullet
#    Do this by inserting well-formed <w:p> paragraphs at the (collapsed)
#    point immediately following the paragraph mark of the target paragraph,
#    i.e. the start of whatever paragraph currently follows it.
# ---------------------------------------------------------------------------
$afterTarget = $d.Paragraphs.Item($targetIndex + 1)
$insertPoint = $d.Range($afterTarget.Range.Start, $afterTarget.Range.Start)

$materialBody =
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Installing angular material for design </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">-  </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> install --save @angular/material @angular/</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>cdk</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> @angular/animations</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr></w:p>'

$materialXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $materialBody + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($materialXml)

# ---------------------------------------------------------------------------
# 3) Move the _GoBack bookmark into its own (already-existing) trailing
#    paragraph - the one that used to directly follow the "Open up" bullet
#    and which is still the last paragraph of the body.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bookmarkPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$bookmarkXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$bookmarkPoint.InsertXML($bookmarkXml)

Write-Output "Edit applied successfully"
